$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3280
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 3725
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 3725
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -4075

$ws.Range("H45").Value = 10000000
$ws.Range("I45").Value = 10000000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 30000000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -29999808
$ws.Range("N45").ClearContents()

$ws.Range("H64").Value = 3921.5208
$ws.Range("I64").Value = 3786.3635
$ws.Range("J64").Value = 4035.8845
$ws.Range("K64").Value = 3786.3635
$ws.Range("L64").Value = 4035.8845
$ws.Range("M64").Value = -3538.3635
$ws.Range("N64").Value = -4531.8845

$ws.Range("H67").Value = 3921.5208
$ws.Range("I67").Value = 3786.3635
$ws.Range("J67").Value = 4035.8845
$ws.Range("K67").Value = 3786.3635
$ws.Range("L67").Value = 4035.8845
$ws.Range("M67").Value = -2928.3635
$ws.Range("N67").Value = -5751.8845

$ws.Range("H86").Value = 1799.1666
$ws.Range("I86").Value = 1798.5714
$ws.Range("K86").Value = 1798.5714
$ws.Range("M86").Value = -675.5714

$ws.Range("H89").Value = 1799.1666
$ws.Range("I89").Value = 1798.5714
$ws.Range("K89").Value = 8992.857
$ws.Range("M89").Value = -3376.857

$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -50120

$ws.Range("H137").Value = 7693556
$ws.Range("I137").Value = 9804931
$ws.Range("J137").Value = 2119
$ws.Range("K137").Value = 29414793
$ws.Range("L137").Value = 6357
$ws.Range("M137").Value = -29412243
$ws.Range("N137").Value = -11457

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2683.5
$ws.Range("I61").Value = 2092.6086
$ws.Range("J61").Value = 5401.6
$ws.Range("K61").Value = 2092.6086
$ws.Range("L61").Value = 5401.6
$ws.Range("M61").Value = -1880.6086
$ws.Range("N61").Value = -5825.6

$ws.Range("H63").Value = 2700
$ws.Range("I63").Value = 2625
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 2625
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -1939
$ws.Range("N63").Value = -4372

$ws.Range("H66").Value = 2700
$ws.Range("I66").Value = 2625
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 13125
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -9693
$ws.Range("N66").Value = -21864

$ws.Range("H136").Value = 2683.5
$ws.Range("I136").Value = 2092.6086
$ws.Range("J136").Value = 5401.6
$ws.Range("K136").Value = 6277.825800000001
$ws.Range("L136").Value = 16204.8
$ws.Range("M136").Value = -3727.825800000001
$ws.Range("N136").Value = -21304.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1678.8723
$ws.Range("I134").Value = 1245.4
$ws.Range("J134").Value = 2943.1667
$ws.Range("K134").Value = 3736.2
$ws.Range("L134").Value = 8829.500100000001
$ws.Range("M134").Value = -1201.2
$ws.Range("N134").Value = -13899.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1811.1111
$ws.Range("I99").Value = 1663.6364
$ws.Range("J99").Value = 2042.8572
$ws.Range("K99").Value = 1663.6364
$ws.Range("L99").Value = 2042.8572
$ws.Range("M99").Value = -165.6364000000001
$ws.Range("N99").Value = -5038.8572

$ws.Range("H126").Value = 1811.1111
$ws.Range("I126").Value = 1663.6364
$ws.Range("J126").Value = 2042.8572
$ws.Range("K126").Value = 4990.9092
$ws.Range("L126").Value = 6128.571599999999
$ws.Range("M126").Value = -2520.9092
$ws.Range("N126").Value = -11068.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2600
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 3400
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 10200
$ws.Range("M80").Value = -2064
$ws.Range("N80").Value = -12072

$ws.Range("H83").Value = 2600
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 3400
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 30600
$ws.Range("M83").Value = -4320
$ws.Range("N83").Value = -39960

$ws.Range("H129").Value = 2081.389
$ws.Range("I129").Value = 1743.3572
$ws.Range("J129").Value = 2296.5
$ws.Range("K129").Value = 5230.071599999999
$ws.Range("L129").Value = 6889.5
$ws.Range("M129").Value = -230.0715999999993
$ws.Range("N129").Value = -16889.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 546.0476
$ws.Range("I97").Value = 546.0476
$ws.Range("K97").Value = 546.0476
$ws.Range("M97").Value = -50.04759999999999

$ws.Range("H101").Value = 19900
$ws.Range("J101").Value = 19900
$ws.Range("L101").Value = 19900
$ws.Range("N101").Value = -26390

$ws.Range("H123").Value = 26058.824
$ws.Range("J123").Value = 26058.824
$ws.Range("L123").Value = 26058.824
$ws.Range("N123").Value = -30958.824

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 47416.332
$ws.Range("J62").Value = 47416.332
$ws.Range("L62").Value = 47416.332
$ws.Range("N62").Value = -48664.332

$ws.Range("H65").Value = 47416.332
$ws.Range("J65").Value = 47416.332
$ws.Range("L65").Value = 142248.996
$ws.Range("N65").Value = -148488.996

$ws.Range("H136").Value = 1553.5111
$ws.Range("I136").Value = 1030.7179
$ws.Range("J136").Value = 4951.6665
$ws.Range("K136").Value = 3092.1537
$ws.Range("L136").Value = 14854.9995
$ws.Range("M136").Value = -542.1537000000003
$ws.Range("N136").Value = -19954.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 19300
$ws.Range("J63").Value = 19300
$ws.Range("L63").Value = 19300
$ws.Range("N63").Value = -20548

$ws.Range("H66").Value = 19300
$ws.Range("J66").Value = 19300
$ws.Range("L66").Value = 57900
$ws.Range("N66").Value = -64140

$ws.Range("H70").Value = 17547.5
$ws.Range("J70").Value = 15000
$ws.Range("L70").Value = 15000
$ws.Range("N70").Value = -15630

$ws.Range("H73").Value = 17547.5
$ws.Range("J73").Value = 15000
$ws.Range("L73").Value = 15000
$ws.Range("N73").Value = -17184

$ws.Range("H132").Value = 2175.0476
$ws.Range("I132").Value = 1076
$ws.Range("J132").Value = 3961
$ws.Range("K132").Value = 3228
$ws.Range("L132").Value = 11883
$ws.Range("M132").Value = -698
$ws.Range("N132").Value = -16943

$ws.Range("H136").Value = 3730.4092
$ws.Range("I136").Value = 3443.5144
$ws.Range("J136").Value = 4846.1113
$ws.Range("K136").Value = 10330.5432
$ws.Range("L136").Value = 14538.3339
$ws.Range("M136").Value = -7780.5432
$ws.Range("N136").Value = -19638.3339
